# Commit: Tue, May 26, 2020 10:05:19 AM
#
# The underlying OOXML diff shows two things changing:
#   1. The table on slide 16 (the "PLENARY - COMPLETE THE MISSING GAPS" slide)
#      switches its table style from the deck's custom style
#      {7B524FFF-73A8-48D0-B5DE-F33275AE8275} to the built-in PowerPoint
#      style {177A0219-4135-40E0-AF75-69244CE583AC}.
#   2. ppt/theme/theme1.xml and ppt/theme/theme2.xml swap contents (a part
#      renumbering artifact with no visible relationship/ref changes).
#      There is no PowerPoint object-model verb that performs a raw
#      part-content swap like that (Master/NotesMaster ApplyTheme and
#      friends only import external .thmx theme files), so only the
#      scriptable table-style change is applied here.

$p = $ppt.ActivePresentation

# Slide 16 -> the graphicFrame holding the 2-column table (shape 3 in
# document order: title, picture, table). Search by HasTable so this
# still finds the right shape even if shape ordering ever shifts.
$s = $p.Slides.Item(16)

$tbl = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tbl = $candidate.Table
    }
}

$tbl.ApplyStyle("{177A0219-4135-40E0-AF75-69244CE583AC}")
